# finish_goods_stocks.xlsx edit:
# - two discontinued products removed from the catalogue
#   (Фп Фиточай "Опалиховский" (БАД) 20х2,0 г / Фп Фиточай "Тибетский" (БАД) 20х2,0  г)
# - stock list re-sorted/refreshed with updated quantities for every remaining item
# - sheet shrinks from 111 to 109 used rows (110/111 deleted), view scrolled/selected accordingly

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the two trailing rows that no longer exist once the two products are removed.
# Deleting bottom-up keeps row numbers of everything above stable while we do it.
$ws.Rows.Item(111).Delete()
$ws.Rows.Item(110).Delete()

# Refreshed name/quantity pairs for rows 2..109 (row 1 header is unchanged).
$data = @(
    @("ВердиоГаст® Растительный комплекс для улучшения пищеварения (БАД ),  капсулы", 80768),
    @("Береза почки 50г", 11239),
    @("Валериана корневища с корнями 50г", 13348),
    @("Эвкалипт прутовидный листья 75г", 16688),
    @("Шиповник плоды низковитаминные 50г", 23405),
    @("Кукуруза столбики с рыльцами 40г", 17689),
    @("Толокнянка листья 50г", 4762),
    @("Чага (березовый гриб) 50г", 19344),
    @("Сб. Фитонефрол (Урологический сбор) 50г", 18016),
    @("Чистотел трава 50г", 12614),
    @("Боярышник плоды 75г", 16962),
    @("Сенна листья 50г", 16141),
    @("Липа цветки 35г", 16344),
    @("Эрва шерстистая трава 30г", 11477),
    @("Сб. Грудной №4 50г", 30296),
    @("Крапива листья 50г", 11261),
    @("Подорожник большой листья 50г", 7840),
    @("Алтей корни 75г", 5264),
    @("Укроп пахучий плоды 50г", 64934),
    @("Шалфей листья 50г", 34636),
    @("Багульник болотный побеги 50г", 13720),
    @("Аир корневища 75г", 7819),
    @("Мята перечная листья 50г", 22522),
    @("Сб. Фитопектол №2 (Грудной сбор №2) 35г", 6294),
    @("Полынь горькая трава 50г", 45932),
    @("Солодка корни 50г", 32831),
    @("Дуба кора 75г", 81175),
    @("Бессмертник песчаный цветки 30г", 27748),
    @("Сб. Фитопектол №1 (Грудной сбор №1) 35г", 5333),
    @("Мать-и-мачеха листья 35г", 29572),
    @("Девясил корневища и корни 50г", 20032),
    @("Ромашка цветки вн 50г", 120284),
    @("Крушина кора 50г", 12984),
    @("Ноготки цветки 50г", 27028),
    @("Лен семена 100г", 67885),
    @("Зверобой трава 50г", 37842),
    @("Пижма цветки 75г", 21055),
    @("Тысячелистник трава 50г", 18079),
    @("Сб. Фитогепатол №2 (Желчегонный сбор №2) 35г", 7041),
    @("Можжевельник плоды 50г", 18422),
    @("Чабрец трава 50г", 28951),
    @("Ламинарии слоевища (морская капуста) 100г", 21681),
    @("Спорыш трава 50г", 23483),
    @("Брусника листья 50г", 24046),
    @("Пустырник трава 50г", 23266),
    @("Череда трава 50г", 25872),
    @("Рябина плоды 50г", 5432),
    @("Фп Детский травяной чай ""ФармаЦветик® для иммунитета"" 20х1,5 г", 2115),
    @("Фп Фиточай ""Лактафитол"" (БАД) 20х1,5 г", 16422),
    @("Фп Детский травяной чай ""ФармаЦветик®  при простуде"" 20х1,5 г", 3745),
    @("Фп Детский травяной чай ""ФармаЦветик® для спокойного сна"" 20х1,5 г", 5400),
    @("Фп Детский травяной чай ""ФармаЦветик® для животика"" 20х1,5 г", 4835),
    @("Фп ""ВердиоГаст® Фиточай для улучшения пищеварения с зеленым чаем""(БАД) 20*1,5г", 6825),
    @("Фп ""ВердиоГаст® Фиточай для улучшения пищеварения с черным чаем"" (БАД) 20*1,5г", 8700),
    @("Фп ""Щедрость природы® Фиточай для иммунитета"" 20х2,0 г", 1602),
    @("Фп ""Щедрость природы® Фиточай кардиологический"" 20х2,0 г", 1602),
    @("Фп ""Щедрость природы® Фиточай при простуде"" 20х2,0 г", 1422),
    @("Фп ""Щедрость природы® Фиточай диабетический"" 20х2,0 г", 1422),
    @("Фп Сб. Грудной №4 20x2,0г", 299813),
    @("Фп ""Щедрость природы® Фиточай успокоительный""20х2,0 г", 2538),
    @("Фп Сб. Бруснивер 20x2,0г", 108757),
    @("Фп Сб. Фитогепатол №3 (Желчегонный сбор №3) 20x2,0г", 50039),
    @("Фп Сенна листья 20x1,5г", 43136),
    @("Фп Сб. Проктофитол (Противогеморроидальный сбор) 20х2,0г", 14272),
    @("Фп Сб. Фитогастрол (Желудочно-кишечный сбор) 20x2,0г", 53469),
    @("Фп Шиповник плоды 20х2,0г", 30384),
    @("Фп Череда трава 20х1,5г", 32363),
    @("Фп Шалфей листья 20х1,5г", 111678),
    @("Фп Ромашка цветки 20x1,5г", 968321),
    @("Фп Мелисса лекарственная трава 20x1,5г", 25506),
    @("Фп Брусника листья 20х1,5г", 55962),
    @("Фп ""Щедрость природы® Фиточай для пищеварения"" 20х2,0 г", 1512),
    @("Фп Чистотел трава 20х1,5г", 23478),
    @("Фп Липа цветки 20x1,5г", 50253),
    @("Фп Чабрец трава 20x1,5 г", 47880),
    @("Фп ""Щедрость природы® Фиточай очищающий"" 20х2,0 г", 2736),
    @("Фп Душица трава 20x1,5г", 20610),
    @("Фп Хвощ полевой трава 20х1,5г", 20698),
    @("Фп Сб. Фитоседан №3 (Успокоительный сбор №3) 20х2,0г", 65517),
    @("Фп Сб. Фитоседан №2 (Успокоительный сбор №2) 20x2,0г", 44695),
    @("Фп Крапива листья 20x1,5г", 52371),
    @("Фп Зверобой трава 20x1,5г", 38957),
    @("Фп Пастушья сумка трава 20х1,5г", 5594),
    @("Фп Толокнянка листья 20x1,5г", 32524),
    @("Фп Мята перечная листья 20x1,5г", 57155),
    @("Фп Сб. Элекасол 20x2,0г", 34068),
    @("Фп Пустырник трава 20x1,5г", 36138),
    @("Фп Сб. Желудочный №3 20x2,0г", 19080),
    @("Фп Сб. Арфазетин-Э 20x2,0г", 44405),
    @("Фп Аир корневища 20x1,5г", 4235),
    @("Фп Сб. Фитонефрол (Урологический сбор) 20x2,0г", 191551),
    @("Фп Золототысячник трава 20х1,5г", 4143),
    @("Фп Фиалка трехцветная трава 20x1,5г", 3870),
    @("Фп Подорожник листья 20x1,5г", 31880),
    @("Фп Береза листья 20x1,5г", 5022),
    @("Фп Боярышник плоды 20х3,0г", 21464),
    @("Фп Пижма цветки 20х1,5г", 9786),
    @("Фп Ольха соплодия 20х1,5г", 4478),
    @("Фп Дуб кора 20х1,5г", 7245),
    @("Фп Ноготки цветки 20x1,5г", 42269),
    @("Фп Валериана корневища с корнями 20x1,5г", 19358),
    @("Фп Лапчатка корневища 20x2,5г", 4063),
    @("Фп Кровохлебка корневища и корни 20x1,5г", 7306),
    @("Фп Крушина кора 20x1,5г", 12383),
    @("Фп Тысячелистник трава 20x1,5г", 24026),
    @("Фп Почечный чай листья 20x1,5г", 111046),
    @("Фп Девясил корневища и корни 20х1,5г", 22874),
    @("Фп Бадан корневища 20x1,5г", 3949)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value2 = $data[$i][0]
    $ws.Cells.Item($row, 2).Value2 = $data[$i][1]
}

# Match the scroll position / active selection recorded after the edit.
try {
    $excel.ActiveWindow.ScrollRow = 64
    $excel.ActiveWindow.ScrollColumn = 1
} catch {}
$ws.Range("A92").Select()

